$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.6052799672350293, 1.7882640872041984, -1.3585632612089646, 2.5423610201704983)"
$ws.Range("C2").Value = "NCT(3.898824241927285, 3.2534165807933397, -4.499285631339608, 3.7886449486043166)"
$ws.Range("D2").Value = "JSU(-0.8582328301860884, 1.0347669633502354, 0.3856610399861583, 2.71149313721872)"
$ws.Range("E2").Value = "NIG(0.766254854856606, 0.5284052195459137, 5.529037748284827, 4.6402262290647025)"
